# "Untitled3333: Examined this question about weighted average wavelength."
#
# The simulation input "b1" (column B of the two result tables, rows 5 and
# 8) was changed from -50 to -10. Every other cell in those two rows is a
# recomputed/derived result of that input, so they are updated to the new
# simulation output values as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First results table (row 4 header / row 5 data): b0, b1, a1, a0, ZaT, ZpT, ZaS, ZpS
$ws.Range("B5").Value = -10
$ws.Range("C5").Value = -63.200008460888483
$ws.Range("D5").Value = 0.55187355689168172
$ws.Range("E5").Value = 15.029296875000005
$ws.Range("F5").Value = 15.000563657221594
$ws.Range("G5").Value = 15.029947916666666
$ws.Range("H5").Value = 15.001203484584281

# Second results table (row 7 header / row 8 data): b0, b1, P - λA/2, ZaS - ZaT, ZpS - ZpT
$ws.Range("B8").Value = -10
$ws.Range("D8").Value = 0.00065104166666074548
$ws.Range("E8").Value = 0.00063982736268641816
